$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark that currently wraps almost the whole document
#    (it will be re-created later, scoped to the new paragraph 6 text).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Merge the "  organization={Springer}" paragraph with the following "}" only
#    paragraph into a single paragraph reading "  organization={Springer}}".
$d.Content.Find.Execute("organization={Springer}" + [char]13 + "}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "organization={Springer}}", 2) | Out-Null

# 3. Insert new clause 6 text into what is now the empty paragraph between
#    clause 5 ("The User will not disclose...") and the empty paragraph that
#    precedes the "Date:" line.
$target = $d.Content.Find
$target.Execute("The User will not disclose the proprietary data set to public or to any third party.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$clausePara = $target.Parent.Paragraphs(1).Next()
$insertRange = $clausePara.Range
$insertRange.Collapse(1)
# A trailing sentinel character is appended after the real clause text so the
# later collapsed-bookmark insertion point does not sit exactly on the
# paragraph-mark boundary (that boundary position is mishandled by
# Bookmarks.Add); the sentinel is stripped again right afterwards.
$insertRange.InsertAfter("6.`tThe User agrees that discloser may revoke the use of all or part of these data at any time. In this case the User is obliged to stop using the data and erase all copies of the data.~")
$sentinelStart = $insertRange.End - 1
$sentinelEnd = $insertRange.End

# 4. Re-create the _GoBack bookmark as a collapsed bookmark right after the
#    newly inserted clause-6 text (before its paragraph mark).
$bmRange = $d.Range($sentinelStart, $sentinelStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 5. Remove the sentinel character now that the bookmark is safely anchored.
$sentinelRange = $d.Range($sentinelStart, $sentinelEnd)
$sentinelRange.Delete()
